# RMA Details Maintenance Grid sheet is refreshed each automation run:
# the three rows of "current" RMA record numbers/ids are overwritten with
# the numbers produced by the latest run (this sync brought in the
# "RMA-391D" run).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RMA Details Maintenance Grid")

$ws.Range("E2").Value = "RMA-391D-001"
$ws.Range("F2").Value = "RMA-391D-1-1"
$ws.Range("J2").Value = "a7s5f000000xLcIAAU"

$ws.Range("E3").Value = "RMA-391D-002"
$ws.Range("F3").Value = "RMA-391D-1-2"
$ws.Range("J3").Value = "a7s5f000000xLcJAAU"

$ws.Range("E4").Value = "RMA-391D-003"
$ws.Range("F4").Value = "RMA-391D-1-3"
$ws.Range("J4").Value = "a7s5f000000xLcKAAU"

# Column widths were re-autofit for the new (slightly different length)
# values in columns F and J.
$ws.Columns.Item(6).ColumnWidth = 13.498697916666666
$ws.Columns.Item(10).ColumnWidth = 19.615885416666668
